$d = $word.ActiveDocument
$d.Content.Find.Execute("how demand vary", $true, $false, $false, $false, $false,
                         $true, 1, $false, "does demand vary", 2)
